# test_data_script output refresh:
#  - collapse the old "Unnamed: 0" index column + numeric running-id column
#  - keep only input_text / output, renamed to columns A / B
#  - re-run the test rows (one new blank/whitespace row inserted after the
#    header) and record PASSED for every non-blank row's output

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Start clean - the whole sheet layout (columns + rows) is being regenerated.
$ws.Cells.Clear()

# Header row
$ws.Range("A1").Value = "input_text"
$ws.Range("B1").Value = "output"

# Data rows (input_text, output) - "" means the output cell stays blank
# (still present in the sheet, just with no text), like the whitespace-only
# input row that never produced a result.
$data = @(
    @(" ", ""),
    @("xin chào", "PASSED"),
    @("https://e-learning.hcmut.edu.vn/course/view.php?id=67808", "PASSED"),
    @("kiemtraphanmem", "PASSED"),
    @("ナム人", "PASSED"),
    @("hello hello", "PASSED"),
    @("xin#@$", "PASSED"),
    @("xin chào", "PASSED"),
    @("a", "PASSED")
)

$r = 2
foreach ($row in $data) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $outCell = $ws.Cells.Item($r, 2)
    if ($row[1] -eq "") {
        # Force the (otherwise empty) cell to materialize without pulling in
        # a new style/number-format, matching the still-present-but-blank
        # output cell from the source data.
        $outCell.Interior.Pattern = 0
    } else {
        $outCell.Value = $row[1]
    }
    $r = $r + 1
}

# Re-apply the bold / bordered / centered header style (matches the shared
# cellXf used by the former header row) - one cell at a time so both land on
# the same final style index.
foreach ($addr in @("A1", "B1")) {
    $c = $ws.Range($addr)
    $c.Font.Bold = $true
    $c.Borders.LineStyle = 1
    $c.HorizontalAlignment = -4108
    $c.VerticalAlignment = -4160
}

Write-Output "ok"
